# Formulas are reliably evaluated and wordwrap and test cleanup
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first")

# Fix typo: "XSLX" -> "XLSX" in C2
$ws.Range("C2").Value = "XLSX"

# New long descriptive cell F15 with word-wrap, and taller row to show it
$ws.Range("F15").Value = "here is some lorem ipsum supporting what should appear as a word wrap. It should NOT be one very long cell."
$ws.Range("F15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 90

# Make sure the formulas below recalc reliably
$ws.Range("C18:C19").Calculate()

# Update the active selection/view (previously scrolled to A2, selecting F15)
$ws.Range("C3").Select()
